$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (iteration result of the linear equation solver)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# Remove the now-unneeded extra iteration rows (3-7)
$ws.Range("A3:D7").ClearContents()
